# Add the "Record_Cooldown" sheet after "Property1" and populate it with the
# cooldown-module drop-item data, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Record_Cooldown"

# --- cell values -----------------------------------------------------
# Order matters: new shared strings are interned in first-seen order, so
# write cells left-to-right, top-to-bottom exactly like Excel would when a
# user types the rows in (Row, Col before Cooldown; SkillID before Time).
$ws2.Range("A1").Value = "Id"
$ws2.Range("A2").Value = "Row"
$ws2.Range("A3").Value = "Col"
$ws2.Range("B1").Value = "Cooldown"
$ws2.Range("B2").Value = 8
$ws2.Range("B3").Value = 2

$ws2.Range("A4").Value = "Public"
$ws2.Range("B4").Value = 0
$ws2.Range("A5").Value = "Private"
$ws2.Range("B5").Value = 1
$ws2.Range("A6").Value = "Save"
$ws2.Range("B6").Value = 0
$ws2.Range("A7").Value = "Cache"
$ws2.Range("B7").Value = 1
$ws2.Range("A8").Value = "Upload"
$ws2.Range("B8").Value = 0

$ws2.Range("A9").Value = "SkillID"
$ws2.Range("B9").Value = "Time"
$ws2.Range("A10").Value = "string"
$ws2.Range("B10").Value = "int"
$ws2.Range("A11").Value = "Desc"

# --- layout ------------------------------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 13.6

# --- selection / active tab ---------------------------------------------
$ws2.Range("D10").Select()
$ws2.Activate()

Write-Host "done"
